# ReviewDoc ingevuld door Carmen
# Fills in the "Klasse" (A) and "gemaakt door" (B) columns for rows 56-60
# on Sheet1, matching the values Carmen added to the review document.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill column A (Klasse) first so the new shared strings are registered
# in the same order as they appear in the source workbook.
$ws.Range("A56").Value = "AssetDAO"
$ws.Range("A57").Value = "JDBCAssetDAO"
$ws.Range("A58").Value = "RootRepository"
$ws.Range("A59").Value = "Asset model"
$ws.Range("A60").Value = "AssetService"

# Then fill column B (gemaakt door), leaving the "Carmen en " variant
# (row 58) for last since it is a new, distinct shared string.
$ws.Range("B56").Value = "Carmen"
$ws.Range("B57").Value = "Carmen"
$ws.Range("B59").Value = "Carmen"
$ws.Range("B60").Value = "Carmen"
$ws.Range("B58").Value = "Carmen en "

# Match the active cell selection left behind in the saved file
$ws.Range("B60").Select()
